$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 1.805874666666667
$ws.Range("H2").Value = 5.417624
$ws.Range("I2").Value = 0.02415265239695089
$ws.Range("J2").Value = 0.02415265239695089
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 7.248785666666667
$ws.Range("N2").Value = 21.746357
$ws.Range("O2").Value = 0.07891374419744837
$ws.Range("P2").Value = 0.07891374419744837
$ws.Range("Q2").Value = 13.09039839952978
$ws.Range("R2").Value = 117.813585595768
$ws.Range("S2").Value = 0.00190597623294287
$ws.Range("T2").Value = 0.00190597623294287
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 1.805874666666667
$ws.Range("H3").Value = 5.417624
$ws.Range("I3").Value = 0.02415265239695089
$ws.Range("J3").Value = 0.02415265239695089
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 37.82684066666667
$ws.Range("N3").Value = 113.480522
$ws.Range("O3").Value = 0.4118010609547572
$ws.Range("P3").Value = 0.4118010609547572
$ws.Range("Q3").Value = 68.31053327996979
$ws.Range("R3").Value = 614.794799519728
$ws.Range("S3").Value = 0.009946087881935834
$ws.Range("T3").Value = 0.009946087881935834
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 1.805874666666667
$ws.Range("H4").Value = 5.417624
$ws.Range("I4").Value = 0.02415265239695089
$ws.Range("J4").Value = 0.02415265239695089
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 38.20927633333334
$ws.Range("N4").Value = 114.627829
$ws.Range("O4").Value = 0.4159644383477588
$ws.Range("P4").Value = 0.4159644383477588
$ws.Range("Q4").Value = 69.00116416203289
$ws.Range("R4").Value = 621.010477458296
$ws.Range("S4").Value = 0.01004664448890633
$ws.Range("T4").Value = 0.01004664448890633
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 1.805874666666667
$ws.Range("H5").Value = 5.417624
$ws.Range("I5").Value = 0.02415265239695089
$ws.Range("J5").Value = 0.02415265239695089
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 8.572171666666666
$ws.Range("N5").Value = 25.716515
$ws.Range("O5").Value = 0.09332075650003555
$ws.Range("P5").Value = 0.09332075650003555
$ws.Range("Q5").Value = 15.48026765115111
$ws.Range("R5").Value = 139.32240886036
$ws.Range("S5").Value = 0.002253943793165853
$ws.Range("T5").Value = 0.002253943793165853
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 53.204531
$ws.Range("H6").Value = 159.613593
$ws.Range("I6").Value = 0.711583459752355
$ws.Range("J6").Value = 0.7115834597523549
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 7.248785666666667
$ws.Range("N6").Value = 21.746357
$ws.Range("O6").Value = 0.07891374419744837
$ws.Range("P6").Value = 0.07891374419744837
$ws.Range("Q6").Value = 385.6682417145224
$ws.Range("R6").Value = 3471.014175430701
$ws.Range("S6").Value = 0.05615371511803265
$ws.Range("T6").Value = 0.05615371511803263
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 53.204531
$ws.Range("H7").Value = 159.613593
$ws.Range("I7").Value = 0.711583459752355
$ws.Range("J7").Value = 0.7115834597523549
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 37.82684066666667
$ws.Range("N7").Value = 113.480522
$ws.Range("O7").Value = 0.4118010609547572
$ws.Range("P7").Value = 0.4118010609547572
$ws.Range("Q7").Value = 2012.559316881727
$ws.Range("R7").Value = 18113.03385193555
$ws.Range("S7").Value = 0.2930308236838766
$ws.Range("T7").Value = 0.2930308236838766
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 53.204531
$ws.Range("H8").Value = 159.613593
$ws.Range("I8").Value = 0.711583459752355
$ws.Range("J8").Value = 0.7115834597523549
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 38.20927633333334
$ws.Range("N8").Value = 114.627829
$ws.Range("O8").Value = 0.4159644383477588
$ws.Range("P8").Value = 0.4159644383477588
$ws.Range("Q8").Value = 2032.9066271644
$ws.Range("R8").Value = 18296.15964447959
$ws.Range("S8").Value = 0.2959934141734434
$ws.Range("T8").Value = 0.2959934141734434
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 53.204531
$ws.Range("H9").Value = 159.613593
$ws.Range("I9").Value = 0.711583459752355
$ws.Range("J9").Value = 0.7115834597523549
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 8.572171666666666
$ws.Range("N9").Value = 25.716515
$ws.Range("O9").Value = 0.09332075650003555
$ws.Range("P9").Value = 0.09332075650003555
$ws.Range("Q9").Value = 456.0783731764882
$ws.Range("R9").Value = 4104.705358588394
$ws.Range("S9").Value = 0.06640550677700237
$ws.Range("T9").Value = 0.06640550677700235
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 1.081716333333333
$ws.Range("H10").Value = 3.245149
$ws.Range("I10").Value = 0.01446740411909589
$ws.Range("J10").Value = 0.01446740411909589
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 7.248785666666667
$ws.Range("N10").Value = 21.746357
$ws.Range("O10").Value = 0.07891374419744837
$ws.Range("P10").Value = 0.07891374419744837
$ws.Range("Q10").Value = 7.841129852465889
$ws.Range("R10").Value = 70.57016867219301
$ws.Range("S10").Value = 0.001141677027855444
$ws.Range("T10").Value = 0.001141677027855444
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 1.081716333333333
$ws.Range("H11").Value = 3.245149
$ws.Range("I11").Value = 0.01446740411909589
$ws.Range("J11").Value = 0.01446740411909589
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 37.82684066666667
$ws.Range("N11").Value = 113.480522
$ws.Range("O11").Value = 0.4118010609547572
$ws.Range("P11").Value = 0.4118010609547572
$ws.Range("Q11").Value = 40.91791138753089
$ws.Range("R11").Value = 368.261202487778
$ws.Range("S11").Value = 0.005957692365504912
$ws.Range("T11").Value = 0.005957692365504913
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 1.081716333333333
$ws.Range("H12").Value = 3.245149
$ws.Range("I12").Value = 0.01446740411909589
$ws.Range("J12").Value = 0.01446740411909589
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 38.20927633333334
$ws.Range("N12").Value = 114.627829
$ws.Range("O12").Value = 0.4159644383477588
$ws.Range("P12").Value = 0.4159644383477588
$ws.Range("Q12").Value = 41.33159829461344
$ws.Range("R12").Value = 371.984384651521
$ws.Range("S12").Value = 0.006017925628749775
$ws.Range("T12").Value = 0.006017925628749776
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 1.081716333333333
$ws.Range("H13").Value = 3.245149
$ws.Range("I13").Value = 0.01446740411909589
$ws.Range("J13").Value = 0.01446740411909589
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 8.572171666666666
$ws.Range("N13").Value = 25.716515
$ws.Range("O13").Value = 0.09332075650003555
$ws.Range("P13").Value = 0.09332075650003555
$ws.Range("Q13").Value = 9.272658103970555
$ws.Range("R13").Value = 83.45392293573499
$ws.Range("S13").Value = 0.001350109096985759
$ws.Range("T13").Value = 0.001350109096985759
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 18.677085
$ws.Range("H14").Value = 56.031255
$ws.Range("I14").Value = 0.2497964837315983
$ws.Range("J14").Value = 0.2497964837315982
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 7.248785666666667
$ws.Range("N14").Value = 21.746357
$ws.Range("O14").Value = 0.07891374419744837
$ws.Range("P14").Value = 0.07891374419744837
$ws.Range("Q14").Value = 135.386186043115
$ws.Range("R14").Value = 1218.475674388035
$ws.Range("S14").Value = 0.01971237581861742
$ws.Range("T14").Value = 0.01971237581861742
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 18.677085
$ws.Range("H15").Value = 56.031255
$ws.Range("I15").Value = 0.2497964837315983
$ws.Range("J15").Value = 0.2497964837315982
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 37.82684066666667
$ws.Range("N15").Value = 113.480522
$ws.Range("O15").Value = 0.4118010609547572
$ws.Range("P15").Value = 0.4118010609547572
$ws.Range("Q15").Value = 706.4951184127901
$ws.Range("R15").Value = 6358.456065715111
$ws.Range("S15").Value = 0.1028664570234399
$ws.Range("T15").Value = 0.1028664570234399
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 18.677085
$ws.Range("H16").Value = 56.031255
$ws.Range("I16").Value = 0.2497964837315983
$ws.Range("J16").Value = 0.2497964837315982
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 38.20927633333334
$ws.Range("N16").Value = 114.627829
$ws.Range("O16").Value = 0.4159644383477588
$ws.Range("P16").Value = 0.4159644383477588
$ws.Range("Q16").Value = 713.6379018661551
$ws.Range("R16").Value = 6422.741116795395
$ws.Range("S16").Value = 0.1039064540566593
$ws.Range("T16").Value = 0.1039064540566593
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 18.677085
$ws.Range("H17").Value = 56.031255
$ws.Range("I17").Value = 0.2497964837315983
$ws.Range("J17").Value = 0.2497964837315982
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 8.572171666666666
$ws.Range("N17").Value = 25.716515
$ws.Range("O17").Value = 0.09332075650003555
$ws.Range("P17").Value = 0.09332075650003555
$ws.Range("Q17").Value = 160.103178852925
$ws.Range("R17").Value = 160.103178852925
$ws.Range("S17").Value = 0.02396973840413007
$ws.Range("T17").Value = 0.02396973840413007